# Hjemme passive tweaks lichtwark deleted values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values changed
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: B2 and D2 values deleted (cells cleared/removed), C2 and E2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -7.3733077377562868
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = -3.0402514088270851

# Row 3 values updated
$ws.Range("B3").Value = -10.616310651571711
$ws.Range("C3").Value = 6.8726977810813397
$ws.Range("D3").Value = -5.1521007428528565
$ws.Range("E3").Value = 22.659194513281399

# Update the selected range to match the new selection extent
$ws.Range("B1:E3").Select()
